$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: AHC30716 / Diploma of Civil Construction Design ----
$ws.Range("A2").Value = "AHC30716"
$ws.Range("B2").Value = "110597F"
$ws.Range("C2").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("D2").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("H2").WrapText = $true
$ws.Range("I2").Value = 16200
$ws.Range("I2").NumberFormat = "#,##0"
$ws.Range("J2").Value = "16,000 tuition fee + 200 handling fee"
$ws.Range("J2").WrapText = $true
$ws.Range("J2").NumberFormat = "#,##0"
$ws.Range("M2").Value = "TAS"
$ws.Rows.Item(2).RowHeight = 45

# ---- Row 3: RII60520 / Advanced Diploma of Civil Construction Design ----
$ws.Range("A3").Value = "RII60520"
$ws.Range("B3").Value = "111826A"
$ws.Range("C3").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("D3").Value = "ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E3").Value = 104
$ws.Range("H3").Value = "88 wks tuition + 16 wks break"
$ws.Range("H3").WrapText = $true
$ws.Range("I3").Value = 25200
$ws.Range("I3").NumberFormat = "#,##0"
$ws.Range("J3").Value = "25,000 tuition fee + 200 handling fee"
$ws.Range("J3").WrapText = $true
$ws.Range("J3").NumberFormat = "#,##0"
$ws.Range("M3").Value = "TAS"
$ws.Rows.Item(3).RowHeight = 45

# ---- Row 4: ICT60220 / Advanced Diploma of Information Technology (Telecommunications Network Engineering) ----
$ws.Range("A4").Value = "ICT60220"
$ws.Range("B4").Value = "111825B"
$ws.Range("C4").Value = "INFORMATION TECHNOLOGY"
$ws.Range("D4").Value = "ADVANCED DIPLOMA OF INFORMATION TECHNOLOGY " + [char]10 + "(TELECOMMUNICATIONS NETWORK ENGINEERING) "
$ws.Range("D4").WrapText = $true
$ws.Range("E4").Value = 104
$ws.Range("H4").Value = "88 wks tuition + 16 wks break"
$ws.Range("H4").WrapText = $true
$ws.Range("I4").Value = 13200
$ws.Range("I4").NumberFormat = "#,##0"
$ws.Range("J4").Value = "13,000 tuition fee + 200 handling fee"
$ws.Range("J4").WrapText = $true
$ws.Range("J4").NumberFormat = "#,##0"
$ws.Range("M4").Value = "TAS"
$ws.Rows.Item(4).RowHeight = 45

# ---- Row 5: RII50520/RII60520 Package ----
$ws.Range("A5").Value = "RII50520/RII60520"
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Value = "111827M/111826A"
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("D5").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN + ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("D5").WrapText = $true
$ws.Range("E5").Value = 104
$ws.Range("H5").Value = "88 wks tuition + 16 wks break"
$ws.Range("H5").WrapText = $true
$ws.Range("I5").Value = 27200
$ws.Range("I5").NumberFormat = "#,##0"
$ws.Range("J5").Value = "27,000 tuition fee + 200 handling fee"
$ws.Range("J5").WrapText = $true
$ws.Range("J5").NumberFormat = "#,##0"
$ws.Range("M5").Value = "TAS"
$ws.Rows.Item(5).RowHeight = 45

# ---- Sheet view / selection update ----
$ws.Range("D18").Select()
